# AddEntrySD.pptx edit script
# 1) Refresh the cached "datetimeFigureOut" footer/date placeholder text from
#    4/7/2018 -> 4/15/2018 across every slide layout, the slide master and
#    the notes master.
# 2) Rename the "ae:AddEvent" labels on the two "AddEvent...Command(Parser)"
#    rectangles on slide 1 to "ae:AddEntry".

$p = $ppt.ActivePresentation

function Update-DatePlaceholderText($container, $oldText, $newText) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
            $isDatePh = $false
        }
        if ($isDatePh -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

# --- Slide master date placeholder ---
Update-DatePlaceholderText $p.SlideMaster "4/7/2018" "4/15/2018"

# --- Every slide layout's date placeholder ---
$layouts = $p.SlideMaster.CustomLayouts
for ($k = 1; $k -le $layouts.Count; $k++) {
    Update-DatePlaceholderText $layouts.Item($k) "4/7/2018" "4/15/2018"
}

# --- Notes master date placeholder (only reachable via HeadersFooters here) ---
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "4/15/2018"

# --- Slide 1: rename the "ae:AddEvent" runs to "ae:AddEntry" ---
# Shape Id 19  ("...Command",       16pt) -> text is split into three runs:
#               "ae" | ":" | "AddEntry" (matches the captured edit exactly).
# Shape Id 82  ("...CommandParser", 12pt) -> simple single-run text swap.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if (-not $shape.HasTextFrame) {
        continue
    }
    $range = $shape.TextFrame.TextRange
    $fullText = $range.Text
    if ($fullText.Length -lt 11 -or $range.Characters(1, 11).Text -ne "ae:AddEvent") {
        continue
    }

    if ($shape.Id -eq 19) {
        # Touch the colon first so the run splits into "ae" | ":" | "AddEvent...",
        # then rewrite the trailing "AddEvent" piece as "AddEntry".
        $colon = $range.Characters(3, 1)
        $colon.Text = ":"
        $tail = $range.Characters(4, 8)
        $tail.Text = "AddEntry"
    } else {
        # Plain in-place replacement of the "ae:AddEvent" run's text.
        $head = $range.Characters(1, 11)
        $head.Text = "ae:AddEntry"
    }
}
